$d = $word.ActiveDocument

# --- Section 1: consolidate the "In ordem / Pré-ordem / Pós-ordem" results ---
# Merge "In ordem" + ":" + "2,3,4,5,7,8,10" runs -> ":2,3,4,5,7,8,10" becomes one run
$d.Content.Find.Execute(":2,3,4,5,7,8,10", $false, $false, $false, $false, $false, $true, 1, $false, ":2,3,4,5,7,8,10", 2)

# Merge "Pré-ordem:" + "7,3,2,4,5,8,10" runs into a single run
$d.Content.Find.Execute("Pré-ordem:7,3,2,4,5,8,10", $false, $false, $false, $false, $false, $true, 1, $false, "Pré-ordem:7,3,2,4,5,8,10", 2)

# Merge "Pós-ordem:" + "2,4,5,3,10,8,7" runs into a single run
$d.Content.Find.Execute("Pós-ordem:2,4,5,3,10,8,7", $false, $false, $false, $false, $false, $true, 1, $false, "Pós-ordem:2,4,5,3,10,8,7", 2)

# Remove the blank paragraph that sits between the "Pré-ordem" and "Pós-ordem" lines
$preOrdemPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Pré-ordem:7,3,2,4,5,8,10")) {
        $preOrdemPara = $p
        break
    }
}
$blank = $preOrdemPara.Next()
$blank.Range.Delete()

# --- Section 2: append the GitHub link block at the end of the document ---
$end = $d.Content.End
$r = $d.Range($end, $end)
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Corpodetexto"/><w:spacing w:before="157" w:line="259" w:lineRule="auto"/><w:ind w:left="104" w:right="521"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Corpodetexto"/><w:spacing w:before="157" w:line="259" w:lineRule="auto"/><w:ind w:left="104" w:right="521"/></w:pPr><w:r><w:t xml:space="preserve">LINK: GITHUB </w:t></w:r><w:r><w:t>https://github.com/BenficaS/EX_EstrutuBinario.git</w:t></w:r></w:p>'
$r.InsertXML($xml)

Write-Output "done"
